$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9: Number=8, task="冲刺后后坐力驱动bug", due="-"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "冲刺后后坐力驱动bug"
$ws.Range("C9").Value = "-"

# Match formatting used by the rest of column A/B/C/D/F (style index 1: center + wrap text)
$ws.Range("A9:C9").HorizontalAlignment = -4108
$ws.Range("A9:C9").WrapText = $true

# Update selection to C10 (matches the new row's "next" cell)
$ws.Range("C10").Select()
